# list.xlsx: replace the old roll-call name list with the new one
# (GTA-style Chinese character names), growing the list from 5 rows
# (1 header + 4 names) to 23 rows (1 header + 22 names), and move the
# active-cell selection down to the new first empty row (A24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> name, for the final A1:A23 roll-call list.
$rowValues = @{
    1  = "名字依次向下"
    2  = "麦克"
    3  = "富兰克林"
    4  = "崔佛"
    5  = "西米恩"
    6  = "莱斯特"
    7  = "史崔奇"
    8  = "拉玛"
    9  = "吉米"
    10 = "戴夫"
    11 = "杰伊"
    12 = "马丁"
    13 = "德温"
    14 = "汤雅"
    15 = "强尼"
    16 = "小罗"
    17 = "陈伟"
    18 = "陈陶"
    19 = "奥尼尔"
    20 = "阿曼达"
    21 = "吉米"
    22 = "大厨"
    23 = "陶艾迪"
}

# Write order matters for shared-string de-dup slot assignment: row 6
# ("莱斯特") must be written before row 5 ("西米恩") so the new unique
# strings land in the table in the same order as the source workbook.
$writeOrder = @(1, 2, 3, 4, 6, 5, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23)

foreach ($row in $writeOrder) {
    $ws.Cells.Item($row, 1).Value = $rowValues[$row]
}

$ws.Range("A24").Select()
